$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Globo"
$ws.Range("B12").Value = "RJ TV 1"
$ws.Range("C12").Value = "Esportes"
$ws.Range("D12").Value = "2025-04-01T11:43"
$ws.Range("E12").Value = "Neutro"
$ws.Range("F12").Value = "Presidente destituído. Após assembleia, Laila Póvoa assume a presidência do Americano. "
